# Update the "harvester" column (B) for rows 3-22 so that it matches the
# value already used in B2 ("H.BROWN") instead of the old placeholder
# "Retrofitted_2385". This makes "Retrofitted_2385" unused, so Excel will
# drop it from the shared strings table automatically on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3:B22").Value = "H.BROWN"

# Update the selected range shown in the saved view to match the edit.
$ws.Range("B3:B22").Select()
